{"js": "// Replace the date string and each three-digit-divided-by-one-digit\n// problem with its updated value. All \"from\" strings are unique within\n// the document, so a body-wide search/replace per pair is safe and\n// order-independent.\nconst replacements = [\n  [\"2025-02-04 Tuesday\", \"2025-02-05 Wednesday\"],\n  [\"125\u00f77=\", \"632\u00f74=\"],\n  [\"888\u00f77=\", \"475\u00f75=\"],\n  [\"359\u00f79=\", \"639\u00f78=\"],\n  [\"601\u00f77=\", \"886\u00f77=\"],\n  [\"968\u00f74=\", \"265\u00f79=\"],\n  [\"566\u00f75=\", \"845\u00f77=\"],\n  [\"684\u00f75=\", \"304\u00f79=\"],\n  [\"176\u00f79=\", \"780\u00f73=\"],\n  [\"405\u00f73=\", \"665\u00f76=\"],\n  [\"701\u00f79=\", \"489\u00f79=\"],\n  [\"237\u00f72=\", \"533\u00f73=\"],\n  [\"966\u00f73=\", \"634\u00f77=\"],\n  [\"221\u00f75=\", \"153\u00f73=\"],\n  [\"711\u00f73=\", \"555\u00f77=\"],\n  [\"396\u00f72=\", \"707\u00f77=\"],\n  [\"483\u00f74=\", \"512\u00f73=\"],\n  [\"357\u00f74=\", \"692\u00f78=\"],\n  [\"318\u00f74=\", \"184\u00f72=\"],\n  [\"404\u00f75=\", \"521\u00f77=\"],\n  [\"622\u00f75=\", \"564\u00f73=\"],\n  [\"175\u00f78=\", \"842\u00f75=\"],\n  [\"477\u00f72=\", \"820\u00f75=\"],\n  [\"832\u00f72=\", \"172\u00f77=\"],\n  [\"910\u00f74=\", \"584\u00f78=\"],\n  [\"764\u00f74=\", \"990\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [from, to] of replacements) {\n  const found = body.search(from, { matchCase: true, matchWholeWord: false });\n  found.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date string and each three-digit-divided-by-one-digit\n# problem with its updated value. All \"from\" strings are unique within\n# the document, so a document-wide Find/Replace per pair is safe and\n# order-independent.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-02-04 Tuesday\", \"2025-02-05 Wednesday\"),\n    @(\"125\u00f77=\", \"632\u00f74=\"),\n    @(\"888\u00f77=\", \"475\u00f75=\"),\n    @(\"359\u00f79=\", \"639\u00f78=\"),\n    @(\"601\u00f77=\", \"886\u00f77=\"),\n    @(\"968\u00f74=\", \"265\u00f79=\"),\n    @(\"566\u00f75=\", \"845\u00f77=\"),\n    @(\"684\u00f75=\", \"304\u00f79=\"),\n    @(\"176\u00f79=\", \"780\u00f73=\"),\n    @(\"405\u00f73=\", \"665\u00f76=\"),\n    @(\"701\u00f79=\", \"489\u00f79=\"),\n    @(\"237\u00f72=\", \"533\u00f73=\"),\n    @(\"966\u00f73=\", \"634\u00f77=\"),\n    @(\"221\u00f75=\", \"153\u00f73=\"),\n    @(\"711\u00f73=\", \"555\u00f77=\"),\n    @(\"396\u00f72=\", \"707\u00f77=\"),\n    @(\"483\u00f74=\", \"512\u00f73=\"),\n    @(\"357\u00f74=\", \"692\u00f78=\"),\n    @(\"318\u00f74=\", \"184\u00f72=\"),\n    @(\"404\u00f75=\", \"521\u00f77=\"),\n    @(\"622\u00f75=\", \"564\u00f73=\"),\n    @(\"175\u00f78=\", \"842\u00f75=\"),\n    @(\"477\u00f72=\", \"820\u00f75=\"),\n    @(\"832\u00f72=\", \"172\u00f77=\"),\n    @(\"910\u00f74=\", \"584\u00f78=\"),\n    @(\"764\u00f74=\", \"990\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $find\n    $rng.Find.Replacement.Text = $replace\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 1\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n    $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
